$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the first sheet.
$ws.Name = "Export as TSV"

# Freeze the header row (row 1) on the first sheet.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Add errorTitle / error message text to the existing data validations.
$ws.Range("I2").Validation.ErrorTitle = "Value must come from list"
$ws.Range("I2").Validation.ErrorMessage = "Value must be one of: imaging."

$ws.Range("J2").Validation.ErrorTitle = "Value must come from list"
$ws.Range("J2").Validation.ErrorMessage = "Value must be one of: seqFISH."

$ws.Range("K2").Validation.ErrorTitle = "Value must come from list"
$ws.Range("K2").Validation.ErrorMessage = "Value must be one of: RNA."

$ws.Range("L2").Validation.ErrorTitle = "Not a boolean"
$ws.Range("L2").Validation.ErrorMessage = 'The values in this column must be "TRUE" or "FALSE".'

$ws.Range("O2").Validation.ErrorTitle = "Not a number"
$ws.Range("O2").Validation.ErrorMessage = "The values in this column must be numbers."

$ws.Range("P2").Validation.ErrorTitle = "Value must come from list"
$ws.Range("P2").Validation.ErrorMessage = "Value must be one of: nm / um."

$ws.Range("Q2").Validation.ErrorTitle = "Not a number"
$ws.Range("Q2").Validation.ErrorMessage = "The values in this column must be numbers."

$ws.Range("R2").Validation.ErrorTitle = "Value must come from list"
$ws.Range("R2").Validation.ErrorMessage = "Value must be one of: nm / um."

$ws.Range("S2").Validation.ErrorTitle = "Not a number"
$ws.Range("S2").Validation.ErrorMessage = "The values in this column must be numbers."

$ws.Range("T2").Validation.ErrorTitle = "Value must come from list"
$ws.Range("T2").Validation.ErrorMessage = "Value must be one of: nm / um."

$ws.Range("W2").Validation.ErrorTitle = "Not an integer"
$ws.Range("W2").Validation.ErrorMessage = "The values in this column must be integers."

$ws.Range("X2").Validation.ErrorTitle = "Not an integer"
$ws.Range("X2").Validation.ErrorMessage = "The values in this column must be integers."

$ws.Range("Y2").Validation.ErrorTitle = "Not an integer"
$ws.Range("Y2").Validation.ErrorMessage = "The values in this column must be integers."

$ws.Range("Z2").Validation.ErrorTitle = "Not an integer"
$ws.Range("Z2").Validation.ErrorMessage = "The values in this column must be integers."

$ws.Range("AA2").Validation.ErrorTitle = "Not an integer"
$ws.Range("AA2").Validation.ErrorMessage = "The values in this column must be integers."

$ws.Range("AB2").Validation.ErrorTitle = "Not an integer"
$ws.Range("AB2").Validation.ErrorMessage = "The values in this column must be integers."
